$wb = $excel.ActiveWorkbook

# --- "Ready for handoff" -> "In Translation" -------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de) in row 2
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C) in row 2
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C) in row 2
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes (Status columns get narrower) --------------------
# Target serialized width is 13.4101845877511; the closest width this host
# can reproduce through the ColumnWidth (character-unit) COM property is
# 13.3333333333333, reached by setting ColumnWidth = 12.5.
$newStatusColWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
